{"js": "// Fix the spacing/punctuation after the DigiKey part number, and rewrite the\n// sentence about the Type A I/O board to say it already exists (rather than\n// being \"worked on\"), that it \"connects\" (not \"will connect\"), and that it\n// has \"similar functionality\" to the Type T board.\n\n// 1) \") .\" -> \").\"\nconst fix1 = context.document.body.search(\") .\", { matchCase: true, matchWholeWord: false });\nfix1.load(\"items\");\nawait context.sync();\nif (fix1.items.length === 0) {\n  throw new Error('Could not find \") .\" to fix.');\n}\nfor (let i = 0; i < fix1.items.length; i++) {\n  fix1.items[i].insertText(\").\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) \"I am working on a Type A I/O board that will connect to the Application\n//    connector at the top of the SYM-1.\" ->\n//    \"I also have a Type A I/O board that connects to the Application\n//    connector at the top of the SYM-1 with similar functionality.\"\nconst oldSentence =\n  \"I am working on a Type A I/O board that will connect to the Application connector at the top of the SYM-1.\";\nconst newSentence =\n  \"I also have a Type A I/O board that connects to the Application connector at the top of the SYM-1 with similar functionality.\";\n\nconst fix2 = context.document.body.search(oldSentence, { matchCase: true, matchWholeWord: false });\nfix2.load(\"items\");\nawait context.sync();\nif (fix2.items.length === 0) {\n  throw new Error(\"Could not find the Type A I/O board sentence to rewrite.\");\n}\nfor (let i = 0; i < fix2.items.length; i++) {\n  fix2.items[i].insertText(newSentence, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix the spacing/punctuation after the DigiKey part number, and rewrite the\n# sentence about the Type A I/O board to say it already exists (rather than\n# being \"worked on\"), that it \"connects\" (not \"will connect\"), and that it\n# has \"similar functionality\" to the Type T board.\n\n$d = $word.ActiveDocument\n\n# 1) \") .\" -> \").\"\n$range1 = $d.Content\n$found1 = $range1.Find.Execute(\") .\", $false, $false, $false, $false, $false, $true, 1, $false, \").\", 2)\nif (-not $found1) {\n    throw 'Could not find \") .\" to fix.'\n}\n\n# 2) \"I am working on a Type A I/O board that will connect to the Application\n#    connector at the top of the SYM-1.\" ->\n#    \"I also have a Type A I/O board that connects to the Application\n#    connector at the top of the SYM-1 with similar functionality.\"\n$oldSentence = \"I am working on a Type A I/O board that will connect to the Application connector at the top of the SYM-1.\"\n$newSentence = \"I also have a Type A I/O board that connects to the Application connector at the top of the SYM-1 with similar functionality.\"\n\n$range2 = $d.Content\n$found2 = $range2.Find.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2)\nif (-not $found2) {\n    throw \"Could not find the Type A I/O board sentence to rewrite.\"\n}\n"}
